$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '34.361.50'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '1.788.67'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '225.97'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = '0.552'
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '32.59'
$ws.Range("E8").Value = '  +1.56%  '
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").Value = '0.0943'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '2.046.24'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.797.77'
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '11.02'
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("E15").Value = '  +1.48%  '
$ws.Range("D16").Value = '34.369.53'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("E17").Value = '  +1.94%  '
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("D20").Value = '243.96'
$ws.Range("E20").Value = '  -0.76%  '
$ws.Range("E21").Value = '  +2.07%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").Value = '165.49'
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("D26").Value = '7.26'
$ws.Range("E26").Value = '  +1.64%  '
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("E30").Value = '  +6.13%  '
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("D35").Value = '2.58'
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("D36").Value = '1.400.98'
$ws.Range("E36").Value = '  -3.29%  '
$ws.Range("D37").Value = '0.675'
$ws.Range("E37").Value = '  +3.22%  '
$ws.Range("E38").Value = '  +2.05%  '
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").Value = '84.43'
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.79'
$ws.Range("E41").Value = '  +2.38%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '2.41'
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("E43").Value = '  +1.91%  '
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").Value = '0.0524'
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("D46").Value = '1.12'
$ws.Range("E46").Value = '  +3.50%  '
$ws.Range("D47").Value = '5.97'
$ws.Range("E47").Value = '  -1.70%  '
$ws.Range("D48").Value = '1.946.81'
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").Value = '104.66'
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("E51").Value = '  -3.35%  '
